# Edit script: split the instructions paragraph into a numbered list item
# that contains a bold/underlined hyperlink-style URL, per the target diff.
$d = $word.ActiveDocument

# Locate the paragraph that introduces the link instructions.
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "נא לפתוח את הקובץ באמצעות דפדפן כלשהו ולאחר מכן:",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target paragraph for edit"
}

# Expand to the whole paragraph (incl. end-of-paragraph mark) so that the
# replacement XML cleanly produces two separate <w:p> elements.
$target = $searchRange.Paragraphs(1).Range

$xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:bidi/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">נא לפתוח את </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>הלינק: ״</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>https://storage.cloud.google.com/xguyor_html_page/check.htm</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>l</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>״</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:bidi/><w:rPr><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> באמצעות דפדפן כלשהו ולאחר מכן:</w:t></w:r></w:p>'

$target.InsertXML($xml)
